$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Value = "Radio"
$ws.Range("S2").Value = "Radio Nova"
$ws.Range("S3").Value = "France Inter"

# Remaining row 1 header cells
$ws.Range("T1").Value = "Journaux"
$ws.Range("U1").Value = "Télé"

$ws.Range("S4").Value = "France bleu"

# Remaining row 2 data cell
$ws.Range("T2").Value = "Le Monde"

# New data cell in existing column K
$ws.Range("K2").Value = "last.fm"

# Bold the new header cells to match existing header row formatting
$ws.Range("S1:U1").Font.Bold = $true

# Column width for new column S (ColumnWidth 10.14 renders as width="11" in OOXML)
$ws.Range("S1").ColumnWidth = 10.14

# Update view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("K2").Select()
